$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 15
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = 6992554
$ws.Cells.Item(15, 3).Value = 'Thailand Premier League'
$ws.Cells.Item(15, 4).Value = 'Thailand Premier League'
$ws.Cells.Item(15, 5).Value = 45158.33333333334
$ws.Cells.Item(15, 6).Value = 'Sukhothai FC'
$ws.Cells.Item(15, 7).Value = 'Trat FC'
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 'D'
$ws.Cells.Item(15, 11).Value = 1.8
$ws.Cells.Item(15, 12).Value = 3.6
$ws.Cells.Item(15, 13).Value = 4.333
$ws.Cells.Item(15, 14).Value = 1.833
$ws.Cells.Item(15, 15).Value = 3.75
$ws.Cells.Item(15, 16).Value = 4
$ws.Cells.Item(15, 17).Value = -0.5
$ws.Cells.Item(15, 18).Value = 1.8
$ws.Cells.Item(15, 19).Value = 2
$ws.Cells.Item(15, 20).Value = 2.75
$ws.Cells.Item(15, 21).Value = 1.975
$ws.Cells.Item(15, 22).Value = 1.825
$ws.Cells.Item(15, 23).Value = -1
$ws.Cells.Item(15, 24).Value = 2.75
$ws.Cells.Item(15, 25).Value = -1
$ws.Cells.Item(15, 26).Value = -1
$ws.Cells.Item(15, 27).Value = 1
$ws.Cells.Item(15, 28).Value = -1
$ws.Cells.Item(15, 29).Value = 0.825

# Row 16
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = 6992550
$ws.Cells.Item(16, 3).Value = 'Thailand Premier League'
$ws.Cells.Item(16, 4).Value = 'Thailand Premier League'
$ws.Cells.Item(16, 5).Value = 45158.33333333334
$ws.Cells.Item(16, 6).Value = 'Buriram United'
$ws.Cells.Item(16, 7).Value = 'Lamphun Warrior FC'
$ws.Cells.Item(16, 8).Value = 3
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 'H'
$ws.Cells.Item(16, 11).Value = 1.166
$ws.Cells.Item(16, 12).Value = 8
$ws.Cells.Item(16, 13).Value = 12
$ws.Cells.Item(16, 14).Value = 1.25
$ws.Cells.Item(16, 15).Value = 6.5
$ws.Cells.Item(16, 16).Value = 8.5
$ws.Cells.Item(16, 17).Value = -1.75
$ws.Cells.Item(16, 18).Value = 1.875
$ws.Cells.Item(16, 19).Value = 1.925
$ws.Cells.Item(16, 20).Value = 3
$ws.Cells.Item(16, 21).Value = 1.975
$ws.Cells.Item(16, 22).Value = 1.825
$ws.Cells.Item(16, 23).Value = 0.25
$ws.Cells.Item(16, 24).Value = -1
$ws.Cells.Item(16, 25).Value = -1
$ws.Cells.Item(16, 26).Value = 0.875
$ws.Cells.Item(16, 27).Value = -1
$ws.Cells.Item(16, 28).Value = 0
$ws.Cells.Item(16, 29).Value = -0

# Row 85
$ws.Cells.Item(85, 1).Value = 83
$ws.Cells.Item(85, 2).Value = 6992623
$ws.Cells.Item(85, 3).Value = 'Thailand Premier League'
$ws.Cells.Item(85, 4).Value = 'Thailand Premier League'
$ws.Cells.Item(85, 5).Value = 45261.375
$ws.Cells.Item(85, 6).Value = 'Ratchaburi FC'
$ws.Cells.Item(85, 7).Value = 'Chiangrai Utd'
$ws.Cells.Item(85, 8).Value = 3
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 'H'
$ws.Cells.Item(85, 11).Value = 1.7
$ws.Cells.Item(85, 12).Value = 3.75
$ws.Cells.Item(85, 13).Value = 4.2
$ws.Cells.Item(85, 14).Value = 1.7
$ws.Cells.Item(85, 15).Value = 3.75
$ws.Cells.Item(85, 16).Value = 4.333
$ws.Cells.Item(85, 17).Value = -0.75
$ws.Cells.Item(85, 18).Value = 1.925
$ws.Cells.Item(85, 19).Value = 1.875
$ws.Cells.Item(85, 20).Value = 2.5
$ws.Cells.Item(85, 21).Value = 1.85
$ws.Cells.Item(85, 22).Value = 1.95
$ws.Cells.Item(85, 23).Value = 0.7
$ws.Cells.Item(85, 24).Value = -1
$ws.Cells.Item(85, 25).Value = -1
$ws.Cells.Item(85, 26).Value = 0.925
$ws.Cells.Item(85, 27).Value = -1
$ws.Cells.Item(85, 28).Value = 0.8500000000000001
$ws.Cells.Item(85, 29).Value = -1

# Row 86
$ws.Cells.Item(86, 1).Value = 84
$ws.Cells.Item(86, 2).Value = 6992620
$ws.Cells.Item(86, 3).Value = 'Thailand Premier League'
$ws.Cells.Item(86, 4).Value = 'Thailand Premier League'
$ws.Cells.Item(86, 5).Value = 45261.375
$ws.Cells.Item(86, 6).Value = 'Uthai Thani FC'
$ws.Cells.Item(86, 7).Value = 'Sukhothai FC'
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 'D'
$ws.Cells.Item(86, 11).Value = 1.95
$ws.Cells.Item(86, 12).Value = 3.5
$ws.Cells.Item(86, 13).Value = 3.4
$ws.Cells.Item(86, 14).Value = 2.1
$ws.Cells.Item(86, 15).Value = 3.4
$ws.Cells.Item(86, 16).Value = 3
$ws.Cells.Item(86, 17).Value = -0.25
$ws.Cells.Item(86, 18).Value = 1.875
$ws.Cells.Item(86, 19).Value = 1.925
$ws.Cells.Item(86, 20).Value = 2.75
$ws.Cells.Item(86, 21).Value = 1.8
$ws.Cells.Item(86, 22).Value = 2
$ws.Cells.Item(86, 23).Value = -1
$ws.Cells.Item(86, 24).Value = 2.4
$ws.Cells.Item(86, 25).Value = -1
$ws.Cells.Item(86, 26).Value = -0.5
$ws.Cells.Item(86, 27).Value = 0.4625
$ws.Cells.Item(86, 28).Value = -1
$ws.Cells.Item(86, 29).Value = 1

# Row 117
$ws.Cells.Item(117, 1).Value = 115
$ws.Cells.Item(117, 2).Value = 7329293
$ws.Cells.Item(117, 3).Value = 'Thailand Premier League'
$ws.Cells.Item(117, 4).Value = 'Thailand Premier League'
$ws.Cells.Item(117, 5).Value = 45288.375
$ws.Cells.Item(117, 6).Value = 'Chonburi'
$ws.Cells.Item(117, 7).Value = 'Bangkok United'
$ws.Cells.Item(117, 8).Value = 0
$ws.Cells.Item(117, 9).Value = 0
$ws.Cells.Item(117, 10).Value = 'D'
$ws.Cells.Item(117, 11).Value = 3.6
$ws.Cells.Item(117, 12).Value = 3.5
$ws.Cells.Item(117, 13).Value = 1.85
$ws.Cells.Item(117, 14).Value = 4.5
$ws.Cells.Item(117, 15).Value = 4
$ws.Cells.Item(117, 16).Value = 1.615
$ws.Cells.Item(117, 17).Value = 0.75
$ws.Cells.Item(117, 18).Value = 1.975
$ws.Cells.Item(117, 19).Value = 1.825
$ws.Cells.Item(117, 20).Value = 3
$ws.Cells.Item(117, 21).Value = 1.85
$ws.Cells.Item(117, 22).Value = 1.95
$ws.Cells.Item(117, 23).Value = -1
$ws.Cells.Item(117, 24).Value = 3
$ws.Cells.Item(117, 25).Value = -1
$ws.Cells.Item(117, 26).Value = 0.9750000000000001
$ws.Cells.Item(117, 27).Value = -1
$ws.Cells.Item(117, 28).Value = -1
$ws.Cells.Item(117, 29).Value = 0.95

# Row 118
$ws.Cells.Item(118, 1).Value = 116
$ws.Cells.Item(118, 2).Value = 7485127
$ws.Cells.Item(118, 3).Value = 'Thailand Premier League'
$ws.Cells.Item(118, 4).Value = 'Thailand Premier League'
$ws.Cells.Item(118, 5).Value = 45288.375
$ws.Cells.Item(118, 6).Value = 'BG Pathum United'
$ws.Cells.Item(118, 7).Value = 'Chiangrai Utd'
$ws.Cells.Item(118, 8).Value = 2
$ws.Cells.Item(118, 9).Value = 2
$ws.Cells.Item(118, 10).Value = 'D'
$ws.Cells.Item(118, 11).Value = 1.5
$ws.Cells.Item(118, 12).Value = 4
$ws.Cells.Item(118, 13).Value = 5.75
$ws.Cells.Item(118, 14).Value = 1.363
$ws.Cells.Item(118, 15).Value = 4.5
$ws.Cells.Item(118, 16).Value = 6.5
$ws.Cells.Item(118, 17).Value = -1.25
$ws.Cells.Item(118, 18).Value = 1.85
$ws.Cells.Item(118, 19).Value = 1.95
$ws.Cells.Item(118, 20).Value = 3
$ws.Cells.Item(118, 21).Value = 1.825
$ws.Cells.Item(118, 22).Value = 1.975
$ws.Cells.Item(118, 23).Value = -1
$ws.Cells.Item(118, 24).Value = 3.5
$ws.Cells.Item(118, 25).Value = -1
$ws.Cells.Item(118, 26).Value = -1
$ws.Cells.Item(118, 27).Value = 0.95
$ws.Cells.Item(118, 28).Value = 0.825
$ws.Cells.Item(118, 29).Value = -1

# Row 167
$ws.Cells.Item(167, 1).Value = 165
$ws.Cells.Item(167, 2).Value = 6992690
$ws.Cells.Item(167, 3).Value = 'Thailand Premier League'
$ws.Cells.Item(167, 4).Value = 'Thailand Premier League'
$ws.Cells.Item(167, 5).Value = 45381.33333333334
$ws.Cells.Item(167, 6).Value = 'Uthai Thani FC'
$ws.Cells.Item(167, 7).Value = 'Prachuap FC'
$ws.Cells.Item(167, 8).Value = 2
$ws.Cells.Item(167, 9).Value = 2
$ws.Cells.Item(167, 10).Value = 'D'
$ws.Cells.Item(167, 11).Value = 2.3
$ws.Cells.Item(167, 12).Value = 3.2
$ws.Cells.Item(167, 13).Value = 2.7
$ws.Cells.Item(167, 14).Value = 2.1
$ws.Cells.Item(167, 15).Value = 3.3
$ws.Cells.Item(167, 16).Value = 3
$ws.Cells.Item(167, 17).Value = -0.25
$ws.Cells.Item(167, 18).Value = 1.85
$ws.Cells.Item(167, 19).Value = 1.95
$ws.Cells.Item(167, 20).Value = 2.75
$ws.Cells.Item(167, 21).Value = 1.925
$ws.Cells.Item(167, 22).Value = 1.875
$ws.Cells.Item(167, 23).Value = -1
$ws.Cells.Item(167, 24).Value = 2.3
$ws.Cells.Item(167, 25).Value = -1
$ws.Cells.Item(167, 26).Value = -0.5
$ws.Cells.Item(167, 27).Value = 0.475
$ws.Cells.Item(167, 28).Value = 0.925
$ws.Cells.Item(167, 29).Value = -1

# Row 168
$ws.Cells.Item(168, 1).Value = 166
$ws.Cells.Item(168, 2).Value = 6992689
$ws.Cells.Item(168, 3).Value = 'Thailand Premier League'
$ws.Cells.Item(168, 4).Value = 'Thailand Premier League'
$ws.Cells.Item(168, 5).Value = 45381.35416666666
$ws.Cells.Item(168, 6).Value = 'Lamphun Warrior FC'
$ws.Cells.Item(168, 7).Value = 'Port FC'
$ws.Cells.Item(168, 8).Value = 2
$ws.Cells.Item(168, 9).Value = 2
$ws.Cells.Item(168, 10).Value = 'D'
$ws.Cells.Item(168, 11).Value = 3.25
$ws.Cells.Item(168, 12).Value = 3.5
$ws.Cells.Item(168, 13).Value = 1.909
$ws.Cells.Item(168, 14).Value = 3
$ws.Cells.Item(168, 15).Value = 3.3
$ws.Cells.Item(168, 16).Value = 2.1
$ws.Cells.Item(168, 17).Value = 0.25
$ws.Cells.Item(168, 18).Value = 1.875
$ws.Cells.Item(168, 19).Value = 1.925
$ws.Cells.Item(168, 20).Value = 3
$ws.Cells.Item(168, 21).Value = 1.925
$ws.Cells.Item(168, 22).Value = 1.875
$ws.Cells.Item(168, 23).Value = -1
$ws.Cells.Item(168, 24).Value = 2.3
$ws.Cells.Item(168, 25).Value = -1
$ws.Cells.Item(168, 26).Value = 0.4375
$ws.Cells.Item(168, 27).Value = -0.5
$ws.Cells.Item(168, 28).Value = 0.925
$ws.Cells.Item(168, 29).Value = -1

# Row 169
$ws.Cells.Item(169, 1).Value = 167
$ws.Cells.Item(169, 2).Value = 6992691
$ws.Cells.Item(169, 3).Value = 'Thailand Premier League'
$ws.Cells.Item(169, 4).Value = 'Thailand Premier League'
$ws.Cells.Item(169, 5).Value = 45381.375
$ws.Cells.Item(169, 6).Value = 'Sukhothai FC'
$ws.Cells.Item(169, 7).Value = 'Ratchaburi FC'
$ws.Cells.Item(169, 8).Value = 1
$ws.Cells.Item(169, 9).Value = 0
$ws.Cells.Item(169, 10).Value = 'H'
$ws.Cells.Item(169, 11).Value = 2.75
$ws.Cells.Item(169, 12).Value = 3.25
$ws.Cells.Item(169, 13).Value = 2.25
$ws.Cells.Item(169, 14).Value = 2.7
$ws.Cells.Item(169, 15).Value = 3.1
$ws.Cells.Item(169, 16).Value = 2.375
$ws.Cells.Item(169, 17).Value = 0
$ws.Cells.Item(169, 18).Value = 2.025
$ws.Cells.Item(169, 19).Value = 1.775
$ws.Cells.Item(169, 20).Value = 2.75
$ws.Cells.Item(169, 21).Value = 1.95
$ws.Cells.Item(169, 22).Value = 1.85
$ws.Cells.Item(169, 23).Value = 1.7
$ws.Cells.Item(169, 24).Value = -1
$ws.Cells.Item(169, 25).Value = -1
$ws.Cells.Item(169, 26).Value = 1.025
$ws.Cells.Item(169, 27).Value = -1
$ws.Cells.Item(169, 28).Value = -1
$ws.Cells.Item(169, 29).Value = 0.8500000000000001

# Row 170
$ws.Cells.Item(170, 1).Value = 168
$ws.Cells.Item(170, 2).Value = 6992694
$ws.Cells.Item(170, 3).Value = 'Thailand Premier League'
$ws.Cells.Item(170, 4).Value = 'Thailand Premier League'
$ws.Cells.Item(170, 5).Value = 45381.41666666666
$ws.Cells.Item(170, 6).Value = 'Trat FC'
$ws.Cells.Item(170, 7).Value = 'BG Pathum United'
$ws.Cells.Item(170, 8).Value = 2
$ws.Cells.Item(170, 9).Value = 1
$ws.Cells.Item(170, 10).Value = 'H'
$ws.Cells.Item(170, 11).Value = 4
$ws.Cells.Item(170, 12).Value = 4
$ws.Cells.Item(170, 13).Value = 1.615
$ws.Cells.Item(170, 14).Value = 3.8
$ws.Cells.Item(170, 15).Value = 4
$ws.Cells.Item(170, 16).Value = 1.65
$ws.Cells.Item(170, 17).Value = 1
$ws.Cells.Item(170, 18).Value = 1.8
$ws.Cells.Item(170, 19).Value = 2
$ws.Cells.Item(170, 20).Value = 3.25
$ws.Cells.Item(170, 21).Value = 2.025
$ws.Cells.Item(170, 22).Value = 1.775
$ws.Cells.Item(170, 23).Value = 2.8
$ws.Cells.Item(170, 24).Value = -1
$ws.Cells.Item(170, 25).Value = -1
$ws.Cells.Item(170, 26).Value = 0.8
$ws.Cells.Item(170, 27).Value = -1
$ws.Cells.Item(170, 28).Value = -0.5
$ws.Cells.Item(170, 29).Value = 0.3875
